$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update current price ("现价") column I for rows 2-42 (2016-10-28 close).
$ws.Range("I2").Value = 3104.27
$ws.Range("I3").Value = 4.4400000000000004
$ws.Range("I4").Value = 3.17
$ws.Range("I5").Value = 17.899999999999999
$ws.Range("I6").Value = 5.65
$ws.Range("I7").Value = 5.26
$ws.Range("I8").Value = 20.04
$ws.Range("I9").Value = 17.920000000000002
$ws.Range("I10").Value = 18.190000000000001
$ws.Range("I11").Value = 16.95
$ws.Range("I12").Value = 16.809999999999999
$ws.Range("I13").Value = 24.75
$ws.Range("I14").Value = 9.17
$ws.Range("I15").Value = 27.32
$ws.Range("I16").Value = 16.79
$ws.Range("I17").Value = 6.95
$ws.Range("I18").Value = 24.91
$ws.Range("I19").Value = 14.23
$ws.Range("I20").Value = 15.48
$ws.Range("I21").Value = 22.3
$ws.Range("I22").Value = 21.55
$ws.Range("I23").Value = 59.24
$ws.Range("I24").Value = 29.91
$ws.Range("I25").Value = 22.92
$ws.Range("I26").Value = 5.13
$ws.Range("I27").Value = 9.41
$ws.Range("I28").Value = 7.36
$ws.Range("I29").Value = 5.18
$ws.Range("I30").Value = 7.09
$ws.Range("I31").Value = 5.23
$ws.Range("I32").Value = 7.38
$ws.Range("I33").Value = 16.48
$ws.Range("I34").Value = 5.01
$ws.Range("I35").Value = 7.17
$ws.Range("I36").Value = 3.02
$ws.Range("I37").Value = 6.8
$ws.Range("I38").Value = 19.309999999999999
$ws.Range("I39").Value = 43.34
$ws.Range("I40").Value = 31.69
$ws.Range("I41").Value = 26.87
$ws.Range("I42").Value = 23.24

# Typing the raw numbers clears the custom currency number format on these
# cells (reverts them to the workbook default "Normal" style), matching the
# source edit exactly.
$ws.Range("I2:I42").Style = "Normal"

# A few rows moved into a different valuation band; re-color J (current
# position) to match by copying the format from a donor cell that already
# carries the target look (PasteSpecial formats only - formula/value is left
# alone and recalculates on its own from the new I values).
$ws.Range("J2").Copy()
$ws.Range("J21").PasteSpecial(-4122)
$ws.Range("J20").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J4").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the view: scrolled back to the top (row 1) at column C, with S9 selected.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("S9").Select()
